# Change_Request_List.xlsx update
# - Bumps the book view window height slightly
# - Row 2 (CRQ ...248819): new NCR date + new NCR number
# - Row 3 (CRQ ...248827): re-purposed into a new "TEST" change record
# - Rows 4-11: old change records removed (contents cleared, row numbering kept)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Change_List")

# --- Row 2 --------------------------------------------------------------
$c = $ws.Range("B2")
$c.NumberFormat = "@"
$c.Value2 = "10-Aug-20"
$ws.Range("K2").Value = "CRQ000000249241"

# --- Row 3 --------------------------------------------------------------
$c = $ws.Range("B3")
$c.NumberFormat = "@"
$c.Value2 = "10-Aug-20"
$ws.Range("C3").Value = "KM Jiaul Islam Jibon"
$ws.Range("D3").Value = "Operational"
$ws.Range("E3").Value = "TEST "
$ws.Range("F3").Value = "MBKLR30,MBKLR27"
$ws.Range("J3").Value = "e.co_Dhaka Metro"
$ws.Range("K3").Value = "CRQ000000249243"
$ws.Range("L3").Value = "Shahriar Mahbub"

# --- Rows 4-11: drop the old entries, keep the row/No. column as-is ----
$ws.Range("B4:L11").ClearContents()
